{"js": "// The paragraph ends with: \"... such as  define and pass-by-reference parameters.\"\n// The edit removes the trailing \" and pass-by-reference parameters.\" text that\n// followed the word \"define\", leaving the sentence ending at \"define\". Word\n// also drops a `_GoBack` bookmark at the point of the last edit when the\n// document is subsequently saved, so we recreate that too.\n\nconst body = context.document.body;\n\n// Mark the last-edit location with Word's automatic \"_GoBack\" bookmark,\n// positioned immediately after the word \"define\" (where the cursor was\n// left once the trailing text was removed).\nconst defineResults = body.search(\"define\", { matchCase: true });\ndefineResults.load(\"text\");\nawait context.sync();\n\nif (defineResults.items.length > 0) {\n  const afterDefine = defineResults.items[0].getRange(\"After\");\n  afterDefine.insertBookmark(\"_GoBack\");\n}\n\n// Remove the trailing \" and pass-by-reference parameters.\" text (the space\n// before \"and\" is a non-breaking space in the document, but Word's search\n// treats it the same as a regular space).\nconst delResults = body.search(\" and pass-by-reference parameters.\", { matchCase: true });\ndelResults.load(\"text\");\nawait context.sync();\n\nif (delResults.items.length > 0) {\n  delResults.items[0].delete();\n}\n\nawait context.sync();\n", "ps1": "# The paragraph ends with: \"... such as  define and pass-by-reference parameters.\"\n# The edit removes the trailing \" and pass-by-reference parameters.\" text that\n# followed the word \"define\", leaving the sentence ending at \"define\". Word\n# also drops a `_GoBack` bookmark at the point of the last edit when the\n# document is subsequently saved, so we recreate that too.\n#\n# NOTE: the bookmark must be added BEFORE the text is deleted - adding it\n# afterwards (once the document's paragraph/range layout has already shifted\n# from the delete) lands it at the wrong spot.\n\n$wdCollapseEnd = 0\n\n$d = $word.ActiveDocument\n\n# Mark the last-edit location with Word's automatic \"_GoBack\" bookmark,\n# positioned immediately after the word \"define\" (where the cursor was\n# left once the trailing text was removed).\n$bmRange = $d.Content\n$found = $bmRange.Find.Execute(\"define\")\nif ($found) {\n    $bmRange.Collapse($wdCollapseEnd)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n\n# Remove the trailing \" and pass-by-reference parameters.\" text (the space\n# before \"and\" is a non-breaking space in the document, but Word's Find\n# treats it the same as a regular space).\n$delRange = $d.Content\n$found2 = $delRange.Find.Execute(\" and pass-by-reference parameters.\")\nif ($found2) {\n    $delRange.Delete()\n}\n"}
